# Auto-applied edits derived from the OOXML diff for Kraken_Profits (8-sheet workbook)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 228.44444
$ws.Range("I2").Value = 247.625
$ws.Range("J2").Value = 75
$ws.Range("K2").Value = 247.625
$ws.Range("L2").Value = 75
$ws.Range("M2").Value = -134.625
$ws.Range("N2").Value = -301
$ws.Range("H32").Value = 10150
$ws.Range("J32").Value = 10150
$ws.Range("L32").Value = 10150
$ws.Range("N32").Value = -10802
$ws.Range("H46").Value = 4500
$ws.Range("I46").Value = 4500
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 13500
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -13381
$ws.Range("H60").Value = 4500
$ws.Range("I60").Value = 4500
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 13500
$ws.Range("L60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -13016
$ws.Range("H113").Value = 3133
$ws.Range("I113").Value = 2799
$ws.Range("K113").Value = 2799
$ws.Range("M113").Value = 455
$ws.Range("H116").Value = 2106
$ws.Range("I116").Value = 2106
$ws.Range("K116").Value = 2106
$ws.Range("M116").Value = 1336
$ws.Range("H130").Value = 99996
$ws.Range("J130").Value = 99996
$ws.Range("L130").Value = 99996
$ws.Range("N130").Value = -110036
$ws.Range("H132").Value = 5175.4
$ws.Range("I132").Value = 4433.154
$ws.Range("K132").Value = 13299.462
$ws.Range("M132").Value = -10769.462
$ws.Range("H138").Value = 3130.72
$ws.Range("I138").Value = 2303.818
$ws.Range("J138").Value = 3780.4285
$ws.Range("K138").Value = 6911.454000000001
$ws.Range("L138").Value = 11341.2855
$ws.Range("M138").Value = -1771.454000000001
$ws.Range("N138").Value = -21621.2855
$ws.Range("H139").Value = 58000
$ws.Range("J139").Value = 58000
$ws.Range("L139").Value = 58000
$ws.Range("N139").Value = -68280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 14999.5
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 14999.5
$ws.Range("K23").Value = 0
$ws.Range("L23").ClearContents()
$ws.Range("M23").Value = 14999.5
$ws.Range("N23").Value = -15517.5
$ws.Range("H55").Value = 16713.715
$ws.Range("I55").Value = 5665.3335
$ws.Range("K55").Value = 5665.3335
$ws.Range("M55").Value = -5350.3335
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").ClearContents()
$ws.Range("N119").Value = 0
$ws.Range("H122").Value = 2858.4614
$ws.Range("I122").Value = 3392
$ws.Range("J122").Value = 2525
$ws.Range("K122").Value = 10176
$ws.Range("L122").Value = 7575
$ws.Range("M122").Value = -7726
$ws.Range("N122").Value = -12475
$ws.Range("H124").Value = 42500
$ws.Range("J124").Value = 42500
$ws.Range("L124").Value = 42500
$ws.Range("N124").Value = -52320
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").ClearContents()
$ws.Range("N134").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 80000
$ws.Range("I57").Value = 80000
$ws.Range("K57").Value = 80000
$ws.Range("M57").Value = -79280
$ws.Range("H59").Value = 80000
$ws.Range("I59").Value = 80000
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 80000
$ws.Range("L59").ClearContents()
$ws.Range("M59").Value = -79153
$ws.Range("N59").Value = 0
$ws.Range("H94").Value = 8836.875
$ws.Range("I94").Value = 8179.2
$ws.Range("K94").Value = 8179.2
$ws.Range("M94").Value = -7728.2
$ws.Range("H99").Value = 10000
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H122").Value = 20780
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 20780
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").Value = 20780
$ws.Range("N122").Value = -30580
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").ClearContents()
$ws.Range("N124").Value = 0
$ws.Range("H129").Value = 75000
$ws.Range("J129").Value = 75000
$ws.Range("L129").Value = 75000
$ws.Range("N129").Value = -85000
$ws.Range("H136").Value = 80000
$ws.Range("I136").Value = 80000
$ws.Range("K136").Value = 80000
$ws.Range("M136").Value = -74900
$ws.Range("H139").Value = 50000
$ws.Range("I139").Value = 50000
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 50000
$ws.Range("L139").ClearContents()
$ws.Range("M139").Value = -44860
$ws.Range("N139").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1400
$ws.Range("I16").Value = 1400
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1400
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -1113
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").ClearContents()
$ws.Range("N20").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").ClearContents()
$ws.Range("N30").Value = 0
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").ClearContents()
$ws.Range("N110").Value = 0
$ws.Range("H113").Value = 1400
$ws.Range("I113").Value = 1400
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1400
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 770
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").ClearContents()
$ws.Range("N128").Value = 0
$ws.Range("H135").Value = 99995
$ws.Range("J135").Value = 99995
$ws.Range("L135").Value = 99995
$ws.Range("N135").Value = -110135
$ws.Range("H137").Value = 80000
$ws.Range("I137").Value = 80000
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 80000
$ws.Range("L137").ClearContents()
$ws.Range("M137").Value = -74900
$ws.Range("N137").Value = 0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1731.6364
$ws.Range("I12").Value = 236.5
$ws.Range("K12").Value = 709.5
$ws.Range("M12").Value = -536.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2829.6
$ws.Range("I97").Value = 3374.5
$ws.Range("J97").Value = 650
$ws.Range("K97").Value = 3374.5
$ws.Range("L97").Value = 650
$ws.Range("M97").Value = -2878.5
$ws.Range("N97").Value = -1642
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").ClearContents()
$ws.Range("N130").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2664
$ws.Range("I7").Value = 2664
$ws.Range("K7").Value = 2664
$ws.Range("M7").Value = -2552
$ws.Range("H46").Value = 3701.5386
$ws.Range("I46").Value = 500
$ws.Range("K46").Value = 500
$ws.Range("M46").Value = -312
$ws.Range("H97").Value = 79672
$ws.Range("J97").Value = 79672
$ws.Range("L97").Value = 79672
$ws.Range("N97").Value = -81654
$ws.Range("H100").Value = 11161.5
$ws.Range("I100").Value = 3823
$ws.Range("K100").Value = 3823
$ws.Range("M100").Value = -3282
$ws.Range("H110").Value = 60000
$ws.Range("J110").Value = 60000
$ws.Range("L110").Value = 60000
$ws.Range("N110").Value = -68180
$ws.Range("H126").Value = 2664
$ws.Range("I126").Value = 2664
$ws.Range("K126").Value = 7992
$ws.Range("M126").Value = -5522
$ws.Range("H134").Value = 20000
$ws.Range("J134").Value = 20000
$ws.Range("L134").Value = 20000
$ws.Range("N134").Value = -30140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("H139").Value = 50000
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
